$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha15")
$ws.Activate()
$ws.Range("A1").Value = "hello"
Write-Host $ws.Name
